$wb = $excel.ActiveWorkbook

# Locate the Norway worksheet (template for the new Italy sheet) and copy it
# to the end of the workbook, right after itself.
$norway = $wb.Worksheets.Item("Norway")
$norway.Copy($null, $norway)

# The copy is created right after "Norway" and named "Norway (2)" -> rename to Italy
$italy = $wb.Worksheets.Item("Norway (2)")
$italy.Name = "Italy"

# Update the market name and the associated NGC ticket reference on the new sheet
$italy.Range("B2").Value = "Italy Market"
$italy.Range("B4").Value = "NGC-3443/T1971/T1927/T1945/T1959"

# The long ticket text needs to wrap, so grow row 4 to fit it
$italy.Range("B4").WrapText = $true
$italy.Rows.Item(4).RowHeight = 43.2

# Put the cell cursor on B4 for the new (now active/selected) Italy sheet
[void]$italy.Range("B4").Select()

# Norway is no longer the active tab; leave it with a "select all" style
# selection, matching a freshly-deselected sheet.
[void]$norway.Range("A1:XFD1048576").Select()

# Make Italy the active sheet/tab
[void]$italy.Activate()
